$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update data values in rows 2-5 (new sensor readings) ---
$ws.Range("A2").Value = 45096.50694444445
$ws.Range("B2").Value = 22.58
$ws.Range("C2").Value = 15.542
$ws.Range("D2").Value = 4.221
$ws.Range("E2").Value = 47.493
$ws.Range("F2").Value = 39.284
$ws.Range("G2").Value = 17.769
$ws.Range("H2").Value = 58.8
$ws.Range("I2").Value = 27.341
$ws.Range("J2").Value = 11.61
$ws.Range("K2").Value = 17.881
$ws.Range("L2").Value = 18.828
$ws.Range("M2").Value = 19.728
$ws.Range("N2").Value = 5.673
$ws.Range("O2").Value = 17.67
$ws.Range("P2").Value = 24.849
$ws.Range("Q2").Value = 14.79
$ws.Range("R2").Value = 3.779
$ws.Range("S2").Value = 2.46
$ws.Range("T2").Value = 261.617
$ws.Range("U2").Value = 49.202
$ws.Range("V2").Value = 16.31
$ws.Range("W2").Value = 32.642
$ws.Range("X2").Value = 17.025
$ws.Range("Y2").Value = 2.109
$ws.Range("Z2").Value = 29.266
$ws.Range("AA2").Value = 14.407
$ws.Range("AB2").Value = 12.944
$ws.Range("AC2").Value = 15.145
$ws.Range("AD2").Value = 19.485
$ws.Range("AE2").Value = 3.64
$ws.Range("AF2").Value = 51.902
$ws.Range("AG2").Value = 9.071
$ws.Range("AH2").Value = 20.391

$ws.Range("A3").Value = 45096.51388888889
$ws.Range("B3").Value = 18.737
$ws.Range("C3").Value = 13.389
$ws.Range("D3").Value = 1.895
$ws.Range("E3").Value = 40.114
$ws.Range("F3").Value = 33.19
$ws.Range("G3").Value = 14.745
$ws.Range("H3").Value = 57.815
$ws.Range("I3").Value = 22.687
$ws.Range("J3").Value = 9.875
$ws.Range("K3").Value = 14.899
$ws.Range("L3").Value = 16.13
$ws.Range("M3").Value = 16.865
$ws.Range("N3").Value = 4.711
$ws.Range("O3").Value = 14.663
$ws.Range("P3").Value = 20.741
$ws.Range("Q3").Value = 12.465
$ws.Range("R3").Value = 1.643
$ws.Range("S3").Value = 1.129
$ws.Range("T3").Value = 215.868
$ws.Range("U3").Value = 41.042
$ws.Range("V3").Value = 13.534
$ws.Range("W3").Value = 27.346
$ws.Range("X3").Value = 14.562
$ws.Range("Y3").Value = 1.765
$ws.Range("Z3").Value = 27.828
$ws.Range("AA3").Value = 11.955
$ws.Range("AB3").Value = 10.749
$ws.Range("AC3").Value = 12.598
$ws.Range("AD3").Value = 16.826
$ws.Range("AE3").Value = 1.294
$ws.Range("AF3").Value = 52.266
$ws.Range("AG3").Value = 7.547
$ws.Range("AH3").Value = 16.921

$ws.Range("A4").Value = 45096.52083333334
$ws.Range("B4").Value = 6.726
$ws.Range("C4").Value = 4.6
$ws.Range("D4").Value = 0.999
$ws.Range("E4").Value = 14.212
$ws.Range("F4").Value = 11.767
$ws.Range("G4").Value = 5.294
$ws.Range("H4").Value = 25.216
$ws.Range("I4").Value = 8.144
$ws.Range("J4").Value = 3.477
$ws.Range("K4").Value = 5.212
$ws.Range("L4").Value = 5.786
$ws.Range("M4").Value = 5.934
$ws.Range("N4").Value = 1.696
$ws.Range("O4").Value = 5.263
$ws.Range("P4").Value = 7.407
$ws.Range("Q4").Value = 4.638
$ws.Range("R4").Value = 1.001
$ws.Range("S4").Value = 0.525
$ws.Range("T4").Value = 72.821
$ws.Range("U4").Value = 14.926
$ws.Range("V4").Value = 4.858
$ws.Range("W4").Value = 9.797000000000001
$ws.Range("X4").Value = 5.309
$ws.Range("Y4").Value = 0.5629999999999999
$ws.Range("Z4").Value = 11.518
$ws.Range("AA4").Value = 4.291
$ws.Range("AB4").Value = 3.958
$ws.Range("AC4").Value = 4.62
$ws.Range("AD4").Value = 5.989
$ws.Range("AE4").Value = 0.784
$ws.Range("AF4").Value = 22.945
$ws.Range("AG4").Value = 2.635
$ws.Range("AH4").Value = 6.075

$ws.Range("A5").Value = 45096.52777777778
$ws.Range("B5").Value = 0.47
$ws.Range("C5").Value = 0.06
$ws.Range("D5").Value = 0.57
$ws.Range("E5").Value = 0.73
$ws.Range("F5").Value = 0.32
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 6.71
$ws.Range("I5").Value = 0.58
$ws.Range("J5").Value = 0.26
$ws.Range("K5").Value = 0.18
$ws.Range("L5").Value = 0.34
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 0.38
$ws.Range("P5").Value = 0.6
$ws.Range("Q5").Value = 0.55
$ws.Range("R5").Value = 0.7
$ws.Range("S5").Value = 0.2
$ws.Range("T5").Value = 0
$ws.Range("U5").Value = 1.49
$ws.Range("V5").Value = 0.35
$ws.Range("W5").Value = 0.96
$ws.Range("X5").Value = 0.61
$ws.Range("Y5").Value = 0.02
$ws.Range("Z5").Value = 2.74
$ws.Range("AA5").Value = 0.31
$ws.Range("AB5").Value = 0.43
$ws.Range("AC5").Value = 0.45
$ws.Range("AD5").Value = 0.33
$ws.Range("AE5").Value = 0.5600000000000001
$ws.Range("AF5").Value = 6.43
$ws.Range("AG5").Value = 0.08
$ws.Range("AH5").Value = 0.45

# --- Delete row 6 (row count reduced from 6 to 5 data points) ---
$ws.Rows.Item(6).Delete()

# --- Update column widths ---
$ws.Columns.Item(2).ColumnWidth = 7.166666666666667
$ws.Columns.Item(3).ColumnWidth = 7.166666666666667
$ws.Columns.Item(5).ColumnWidth = 7.166666666666667
$ws.Columns.Item(6).ColumnWidth = 7.166666666666667
$ws.Columns.Item(7).ColumnWidth = 7.166666666666667
$ws.Columns.Item(8).ColumnWidth = 7.166666666666667
$ws.Columns.Item(9).ColumnWidth = 7.166666666666667
$ws.Columns.Item(11).ColumnWidth = 7.166666666666667
$ws.Columns.Item(12).ColumnWidth = 7.166666666666667
$ws.Columns.Item(13).ColumnWidth = 7.166666666666667
$ws.Columns.Item(15).ColumnWidth = 7.166666666666667
$ws.Columns.Item(16).ColumnWidth = 7.166666666666667
$ws.Columns.Item(17).ColumnWidth = 7.166666666666667
$ws.Columns.Item(19).ColumnWidth = 6.166666666666667
$ws.Columns.Item(20).ColumnWidth = 8.166666666666666
$ws.Columns.Item(21).ColumnWidth = 7.166666666666667
$ws.Columns.Item(22).ColumnWidth = 7.166666666666667
$ws.Columns.Item(23).ColumnWidth = 7.166666666666667
$ws.Columns.Item(24).ColumnWidth = 7.166666666666667
$ws.Columns.Item(26).ColumnWidth = 7.166666666666667
$ws.Columns.Item(27).ColumnWidth = 7.166666666666667
$ws.Columns.Item(28).ColumnWidth = 7.166666666666667
$ws.Columns.Item(29).ColumnWidth = 7.166666666666667
$ws.Columns.Item(30).ColumnWidth = 7.166666666666667
$ws.Columns.Item(32).ColumnWidth = 7.166666666666667
$ws.Columns.Item(34).ColumnWidth = 7.166666666666667
